$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format before writing so numeric-looking
# strings (e.g. "178.32", "0.705") are preserved verbatim instead of being
# coerced into floating point numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.818.37"
$ws.Range("E2").Value = "  -6.08%  "
$ws.Range("D3").Value = "3.690.80"
$ws.Range("E3").Value = "  -5.50%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "571.12"
$ws.Range("E5").Value = "  -4.33%  "
$ws.Range("D6").Value = "178.32"
$ws.Range("E6").Value = "  +5.72%  "
$ws.Range("D7").Value = "3.691.87"
$ws.Range("E7").Value = "  -5.05%  "
$ws.Range("D8").Value = "0.624"
$ws.Range("E8").Value = "  -7.95%  "
$ws.Range("D9").Value = "0.998"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").Value = "0.705"
$ws.Range("E10").Value = "  -8.61%  "
$ws.Range("E11").Value = "  -13.58%  "
$ws.Range("D12").Value = "51.76"
$ws.Range("E12").Value = "  -5.56%  "
$ws.Range("E13").Value = "  -12.10%  "
$ws.Range("D14").Value = "10.40"
$ws.Range("E14").Value = "  -8.80%  "
$ws.Range("D15").Value = "4.271.85"
$ws.Range("E15").Value = "  -5.63%  "
$ws.Range("D16").Value = "3.684.68"
$ws.Range("E16").Value = "  -5.88%  "
$ws.Range("E17").Value = "  -8.81%  "
$ws.Range("E18").Value = "  -3.16%  "
$ws.Range("D19").Value = "12.77"
$ws.Range("E19").Value = "  -8.92%  "
$ws.Range("E20").Value = "  -8.49%  "
$ws.Range("D21").Value = "67.445.59"
$ws.Range("E21").Value = "  -6.64%  "
$ws.Range("D22").Value = "405.12"
$ws.Range("E22").Value = "  -8.54%  "
$ws.Range("D23").Value = "4.41"
$ws.Range("E23").Value = "  -6.94%  "
$ws.Range("D24").Value = "87.82"
$ws.Range("E24").Value = "  -7.32%  "
$ws.Range("D25").Value = "3.03"
$ws.Range("E25").Value = "  -8.10%  "
$ws.Range("D26").Value = "12.65"
$ws.Range("E26").Value = "  -9.33%  "
$ws.Range("D27").Value = "10.78"
$ws.Range("E27").Value = "  -3.09%  "
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("E29").Value = "  -9.58%  "
$ws.Range("D30").Value = "9.46"
$ws.Range("E30").Value = "  -7.49%  "
$ws.Range("D31").Value = "32.39"
$ws.Range("E31").Value = "  -8.59%  "
$ws.Range("D32").Value = "7.42"
$ws.Range("E32").Value = "  -5.43%  "
$ws.Range("D33").Value = "12.38"
$ws.Range("E33").Value = "  -9.52%  "
$ws.Range("D34").Value = "610.60"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.115"
$ws.Range("E35").Value = "  -8.99%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "64.66"
$ws.Range("E36").Value = "  -5.24%  "
$ws.Range("D37").Value = "43.00"
$ws.Range("E37").Value = "  -15.50%  "
$ws.Range("E38").Value = "  -9.94%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").Value = "0.394"
$ws.Range("E40").Value = "  -6.71%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "0.134"
$ws.Range("E42").Value = "  -6.39%  "
$ws.Range("D43").Value = "2.75"
$ws.Range("E43").Value = "  +7.22%  "
$ws.Range("D44").Value = "2.98"
$ws.Range("E44").Value = "  -10.22%  "
$ws.Range("E45").Value = "  -8.90%  "
$ws.Range("E46").Value = "  -11.58%  "
$ws.Range("D47").Value = "9.17"
$ws.Range("E47").Value = "  -11.63%  "
$ws.Range("D48").Value = "2.792.15"
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").Value = "0.133"
$ws.Range("E49").Value = "  -8.82%  "
$ws.Range("E50").Value = "  -7.69%  "
$ws.Range("D51").Value = "3.08"
$ws.Range("E51").Value = "  -7.54%  "

# Restore default styling on the touched range so no stray cell format
# lingers from the temporary Text NumberFormat above.
$ws.Range("D2:E51").Style = "Normal"
